$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update a few existing text values (F2 / G2 stay untouched) ---
$ws.Range("B2").Value = "егорdsa2"
$ws.Range("C2").Value = "губин2"
$ws.Range("D2").Value = "выфывфы2"
$ws.Range("E2").Value = "выфв2"

# --- Row 3: replace all text values ---
$ws.Range("B3").Value = "егор1"
$ws.Range("C3").Value = "губин1"
# D3:G3 are purely-numeric-looking text, so lead with an apostrophe to
# force Excel to keep them as text rather than silently parsing numbers.
$ws.Range("D3").Value = "'21"
$ws.Range("E3").Value = "'21"
$ws.Range("F3").Value = "'21"
$ws.Range("G3").Value = "'21"

# --- Row 4: replace all text values (all numeric-looking) ---
$ws.Range("B4").Value = "'3"
$ws.Range("C4").Value = "'3"
$ws.Range("D4").Value = "'3"
$ws.Range("E4").Value = "'3"
$ws.Range("F4").Value = "'3"
$ws.Range("G4").Value = "'3"

# --- New rows 5-7: same index-column style as the existing data rows ---
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'14"
$ws.Range("C5").Value = "'14"
$ws.Range("D5").Value = "'14"
$ws.Range("E5").Value = "'14"
$ws.Range("F5").Value = "'41"
$ws.Range("G5").Value = "'14"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'1"
$ws.Range("C6").Value = "'1"
$ws.Range("D6").Value = "'1"
$ws.Range("E6").Value = "'1"
$ws.Range("F6").Value = "'1"
$ws.Range("G6").Value = "'1"

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'36"
$ws.Range("C7").Value = "'36"
$ws.Range("D7").Value = "'36"
$ws.Range("E7").Value = "'36"
$ws.Range("F7").Value = "'36"
$ws.Range("G7").Value = "'36"
